# Applies the changes described by the commit diff:
#  - Corrects the swapped "Data_Inicio"/"Data_Fim" column headers (G1/H1),
#    relabeling them as "Data Fim" / "Data Início", and renames the
#    "Tempo_Parada_h" header (J1) to "Tempo de Parada (h)".
#  - Clears the two stray duration formulas (J16, J20) that were computing
#    garbage because their matching "Data_Fim"/end-date cell (H16/H20) is
#    blank.
#  - Converts the data range into a native Excel Table ("Tabela1") with an
#    AutoFilter, mirroring the table that ships with the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header labels -----------------------------------------------------
$ws.Range("G1").Value = "Data Fim"
$ws.Range("H1").Value = "Data Início"
$ws.Range("J1").Value = "Tempo de Parada (h)"

# --- Clear the bogus duration formulas on rows with no end date ------------
$ws.Range("J16").ClearContents()
$ws.Range("J20").ClearContents()

# --- Turn the data range into a proper Excel Table --------------------------
$dataRange = $ws.Range("A1:K21")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Tabela1"
$tbl.TableStyle = "TableStyleLight11"

Write-Host "done"
